# CHG: finish homeworktable when newuser give required course AND SQL required column change to 0,1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("功能表")

# The "功課表" (homework/schedule table) row now has its completed-status
# changed from "X" to "O" in column B (row 12).
$ws.Range("B12").Value = "O"

# Update the active selection to reflect where the user ended up working.
$ws.Range("B20").Select()
